$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The shared-string table gained a new entry "Conserved stretch 12" (inserted
# just before the former "Conserved stretch 13"), lost the entry
# "Conserved stretch 120" (it was a stray/duplicate tail entry), and the
# label "ORF1ab, S" was reworded to "S, ORF1ab". Because the worksheet's
# cells keep pointing at the same (now shifted) shared-string slots, the
# net visible effect is that every "Conserved stretch N" label in column A
# (rows 2-109) drops by one, and the gene label in column E for the stretch
# that used to read "ORF1ab, S" now reads "S, ORF1ab".

for ($row = 2; $row -le 109; $row++) {
    $n = $row + 10
    $ws.Cells.Item($row, 1).Value = "Conserved stretch $n"
}

$ws.Range("E79").Value = "S, ORF1ab"
